# TC01_CDS_phs001287_Gender_Female.xlsx — "1437 and 1287 working versions"
#
# Semantic changes applied:
#  1. The "File Name" query (cell B4 on Sheet1) gains a new trailing
#     output column:  '' AS "Supplementary Files"   (added right after
#     gi.library_strategy AS "Library Strategy", just before the FROM clause).
#  2. The saved view resets: zoom goes from 140% back to 100% and the
#     window scrolls back so A1 is the top-left cell again (instead of A3).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Update the "File Name" query text (row 4, column B) -------------
$oldFragment = 'gi.library_strategy AS "Library Strategy"' + [char]10 + 'FROM '
$newFragment = 'gi.library_strategy AS "Library Strategy",' + [char]10 + '    ''' + ''' AS "Supplementary Files"' + [char]10 + 'FROM '

$cell = $ws.Cells.Item(4, 2)
$current = $cell.Value()
if ($current.Contains($oldFragment)) {
    $cell.Value = $current.Replace($oldFragment, $newFragment)
}

# Row auto-sizes when the text grows; keep it pinned at the same
# (already-maxed-out) wrapped height used before the edit.
$ws.Rows.Item(4).RowHeight = 409.6

# --- 2. Reset the window view: zoom 140% -> 100%, scroll back to A1 -----
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$win.Zoom = 100
